# "Worked on allowing MF6 Time Series in Edit Feature Formula dialog box."
#
# Bumps the "Actual (hours)" figure logged against the "Time Lists" row
# (row 3) from 51 to 56 to reflect the extra time spent on MF6 time-series
# support. The downstream rollups in column I (I39:I42 -- total, days,
# months, final days) are formulas that depend on this cell, so they
# recalculate automatically once the new value is entered.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I3").Value = 56

$excel.Calculate()
